$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 824.64703
$ws.Range("I19").Value = 717
$ws.Range("J19").Value = 978.4286
$ws.Range("K19").Value = 717
$ws.Range("L19").Value = 978.4286
$ws.Range("M19").Value = -542
$ws.Range("N19").Value = -1328.4286

$ws.Range("H103").Value = 1382.6666
$ws.Range("I103").Value = 1692
$ws.Range("J103").Value = 300
$ws.Range("K103").Value = 5076
$ws.Range("L103").Value = 900
$ws.Range("M103").Value = -4490
$ws.Range("N103").Value = -2072

$ws.Range("H138").Value = 1731.2632
$ws.Range("I138").Value = 2181.5
$ws.Range("J138").Value = 1523.4615
$ws.Range("K138").Value = 6544.5
$ws.Range("L138").Value = 4570.3845
$ws.Range("M138").Value = -1404.5
$ws.Range("N138").Value = -14850.3845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 59507.824
$ws.Range("I2").Value = 710.9091
$ws.Range("J2").Value = 167302.17
$ws.Range("K2").Value = 710.9091
$ws.Range("L2").Value = 167302.17
$ws.Range("M2").Value = -597.9091
$ws.Range("N2").Value = -167528.17

$ws.Range("H28").Value = 12760.375
$ws.Range("I28").Value = 12760.375
$ws.Range("K28").Value = 12760.375
$ws.Range("M28").Value = -12568.375

$ws.Range("H46").Value = 3500
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = -2681
$ws.Range("N46").Value = -4638

$ws.Range("H74").Value = 1427.55
$ws.Range("I74").Value = 1536.7333
$ws.Range("J74").Value = 1100
$ws.Range("K74").Value = 1536.7333
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = -662.7333000000001
$ws.Range("N74").Value = -2848

$ws.Range("H77").Value = 1427.55
$ws.Range("I77").Value = 1536.7333
$ws.Range("J77").Value = 1100
$ws.Range("K77").Value = 7683.6665
$ws.Range("L77").Value = 5500
$ws.Range("M77").Value = -3315.6665
$ws.Range("N77").Value = -14236

$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0

$ws.Range("H96").Value = 34400
$ws.Range("J96").Value = 34400
$ws.Range("L96").Value = 34400
$ws.Range("N96").Value = -39892

$ws.Range("H97").Value = 28819.5
$ws.Range("I97").Value = 37748.184
$ws.Range("J97").Value = 2033.4445
$ws.Range("K97").Value = 37748.184
$ws.Range("L97").Value = 2033.4445
$ws.Range("M97").Value = -37252.184
$ws.Range("N97").Value = -3025.4445

$ws.Range("H98").Value = 19995
$ws.Range("J98").Value = 19995
$ws.Range("L98").Value = 19995
$ws.Range("N98").Value = -25985

$ws.Range("H99").Value = 12760.375
$ws.Range("I99").Value = 12760.375
$ws.Range("K99").Value = 12760.375
$ws.Range("M99").Value = -9765.375

$ws.Range("H103").Value = 31500
$ws.Range("J103").Value = 31500
$ws.Range("L103").Value = 31500
$ws.Range("N103").Value = -33844

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").ClearContents()
$ws.Range("N104").Value = 0

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").ClearContents()
$ws.Range("N106").Value = 0

$ws.Range("H116").Value = 59507.824
$ws.Range("I116").Value = 710.9091
$ws.Range("J116").Value = 167302.17
$ws.Range("K116").Value = 710.9091
$ws.Range("L116").Value = 167302.17
$ws.Range("M116").Value = 1583.0909
$ws.Range("N116").Value = -171890.17

$ws.Range("H122").Value = 2149.9333
$ws.Range("I122").Value = 1980.75
$ws.Range("J122").Value = 2826.6667
$ws.Range("K122").Value = 5942.25
$ws.Range("L122").Value = 8480.000100000001
$ws.Range("M122").Value = -3492.25
$ws.Range("N122").Value = -13380.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 59507.824
$ws.Range("I3").Value = 710.9091
$ws.Range("J3").Value = 167302.17
$ws.Range("K3").Value = 710.9091
$ws.Range("L3").Value = 167302.17
$ws.Range("M3").Value = -596.9091
$ws.Range("N3").Value = -167530.17

$ws.Range("H80").Value = 1099.1538
$ws.Range("I80").Value = 750
$ws.Range("J80").Value = 1227.7894
$ws.Range("K80").Value = 750
$ws.Range("L80").Value = 1227.7894
$ws.Range("M80").Value = 248
$ws.Range("N80").Value = -3223.7894

$ws.Range("H83").Value = 1099.1538
$ws.Range("I83").Value = 750
$ws.Range("J83").Value = 1227.7894
$ws.Range("K83").Value = 3750
$ws.Range("L83").Value = 6138.946999999999
$ws.Range("M83").Value = 1242
$ws.Range("N83").Value = -16122.947

$ws.Range("H107").Value = 66698588
$ws.Range("I107").Value = 111162620
$ws.Range("K107").Value = 111162620
$ws.Range("M107").Value = -111160700

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1492.7858
$ws.Range("J33").Value = 2948.5715
$ws.Range("L33").Value = 17691.429
$ws.Range("N33").Value = -18257.429

$ws.Range("H34").Value = 1366.5
$ws.Range("J34").Value = 1366.5
$ws.Range("L34").Value = 4099.5
$ws.Range("N34").Value = -4267.5

$ws.Range("H131").Value = 758.52
$ws.Range("J131").Value = 789.5
$ws.Range("L131").Value = 2368.5
$ws.Range("N131").Value = -12448.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2546.5356
$ws.Range("I102").Value = 1821.6666
$ws.Range("K102").Value = 1821.6666
$ws.Range("M102").Value = -199.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2676.6924
$ws.Range("I22").Value = 2599.5715
$ws.Range("J22").Value = 2766.6667
$ws.Range("K22").Value = 2599.5715
$ws.Range("L22").Value = 2766.6667
$ws.Range("M22").Value = -2304.5715
$ws.Range("N22").Value = -3356.6667

$ws.Range("H27").Value = 2676.6924
$ws.Range("I27").Value = 2599.5715
$ws.Range("J27").Value = 2766.6667
$ws.Range("K27").Value = 2599.5715
$ws.Range("L27").Value = 2766.6667
$ws.Range("M27").Value = -2492.5715
$ws.Range("N27").Value = -2980.6667

$ws.Range("H55").Value = 335020.06
$ws.Range("I55").Value = 669452.6
$ws.Range("J55").Value = 587.4706
$ws.Range("K55").Value = 669452.6
$ws.Range("L55").Value = 587.4706
$ws.Range("M55").Value = -669279.6
$ws.Range("N55").Value = -933.4706

$ws.Range("H122").Value = 2700
$ws.Range("I122").Value = 2700
$ws.Range("K122").Value = 8100
$ws.Range("M122").Value = -5650

$ws.Range("H132").Value = 4622.222
$ws.Range("I132").Value = 4622.222
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13866.666
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -11336.666

$ws.Range("H136").Value = 2517.7058
$ws.Range("I136").Value = 2335.7856
$ws.Range("J136").Value = 3366.6667
$ws.Range("K136").Value = 7007.3568
$ws.Range("L136").Value = 10100.0001
$ws.Range("M136").Value = -4457.3568
$ws.Range("N136").Value = -15200.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2629.3235
$ws.Range("I132").Value = 2593.9333
$ws.Range("J132").Value = 2894.75
$ws.Range("K132").Value = 7781.7999
$ws.Range("L132").Value = 8684.25
$ws.Range("M132").Value = -5251.7999
$ws.Range("N132").Value = -13744.25
